# Update gh-pages to output generated at 456a3b4
# Apply the same set of value updates to both the "展览" sheet and the
# "全部类型" sheet (the latter aggregates all rows, so the same events
# show up on different row numbers).

$wb = $excel.ActiveWorkbook

function Update-ConventionSheet {
    param($ws, $rowF2, $rowF4, $rowF5, $rowF6, $rowF7, $rowF8)

    $ws.Cells.Item($rowF2, 6).Value = 13760    # F: 13701 -> 13760
    $ws.Cells.Item($rowF4, 6).Value = 665      # F: 663 -> 665
    $ws.Cells.Item($rowF5, 6).Value = 233      # F: 230 -> 233
    $ws.Cells.Item($rowF5, 7).Value = "不可售"  # G: 40 -> "不可售" (text)
    $ws.Cells.Item($rowF6, 6).Value = 490      # F: 481 -> 490
    $ws.Cells.Item($rowF7, 6).Value = 1404     # F: 1399 -> 1404
    $ws.Cells.Item($rowF8, 6).Value = 134      # F: 131 -> 134
}

# "展览" sheet: rows 2, 4, 5, 6, 7, 8
$wsExhibit = $wb.Worksheets.Item("展览")
Update-ConventionSheet $wsExhibit 2 4 5 6 7 8

# "全部类型" sheet: same events appear at rows 2, 4, 5, 8, 9, 11
$wsAll = $wb.Worksheets.Item("全部类型")
Update-ConventionSheet $wsAll 2 4 5 8 9 11
